# https配置说明.docx — apply the 20170221 commit:
#   - mark several ascii tokens (root.keystore / cmd / keytool / keystore /
#     storepass / changeit / jdk / cmd) as spell-check "not in dictionary"
#     runs (w:proofErr spellStart/spellEnd, plus gramStart/gramEnd around
#     the first "keytool" token) by splitting the surrounding runs:
#   - relocate the _GoBack bookmark from the "https://localhost:8443/"
#     paragraph to the "结束" paragraph.

$d = $word.ActiveDocument

function Insert-Fragment($Range, $InnerXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $InnerXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg)
}

# Common run-properties blocks re-used below.
$rPrHint   = '<w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>'
$rPrMono   = '<w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:cs="宋体"/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="13"/><w:szCs w:val="13"/></w:rPr>'
$rPrMonoHint = '<w:rPr><w:rFonts w:ascii="宋体" w:eastAsia="宋体" w:cs="宋体" w:hint="eastAsia"/><w:bCs/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="13"/><w:szCs w:val="13"/></w:rPr>'

# ---------------------------------------------------------------------
# 1) "root.keystore" bullet -> wrap the run in spellStart/spellEnd
# ---------------------------------------------------------------------
$p = $d.Paragraphs(2)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrHint + '<w:t>root.keystore</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Insert-Fragment $r $inner

# ---------------------------------------------------------------------
# 2) "cmd执行命令导入根证书" bullet -> wrap just the "cmd" run
# ---------------------------------------------------------------------
$p = $d.Paragraphs(7)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrHint + '<w:t>cmd</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPrHint + '<w:t>执行命令导入根证书</w:t></w:r>'
Insert-Fragment $r $inner

# ---------------------------------------------------------------------
# 3) keytool command line -> split into annotated runs
# ---------------------------------------------------------------------
$p = $d.Paragraphs(8)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r>' + $rPrMono + '<w:t>keytool</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r>' + $rPrMono + '<w:t xml:space="preserve"> -import -alias root -file root.crt -</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrMono + '<w:t>keystore</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPrMono + '<w:t xml:space="preserve"> F:\jdk1.6.0_06\jre\lib\security\cacerts -</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrMono + '<w:t>storepass</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPrMono + '<w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrMono + '<w:t>changeit</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Insert-Fragment $r $inner

# ---------------------------------------------------------------------
# 4) "注意：……" / "（在cmd默认执行目录下可不加路径）" paragraph -> split
#    into annotated runs (jdk / cmd marked as spelling exceptions)
# ---------------------------------------------------------------------
$p = $d.Paragraphs(9)
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$inner = '<w:r>' + $rPrMonoHint + '<w:t>注意：用自己的</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrMonoHint + '<w:t>jdk</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPrMonoHint + '<w:t>目录(安装环境不要有空格如Program Files),root.crt也是有路径的</w:t></w:r>' +
         '<w:r>' + $rPrMonoHint + '<w:t>（在</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r>' + $rPrMonoHint + '<w:t>cmd</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r>' + $rPrMonoHint + '<w:t>默认执行目录下可不加路径）</w:t></w:r>'
Insert-Fragment $r $inner

# ---------------------------------------------------------------------
# 5) Move the _GoBack bookmark from the "访问https://localhost:8443/"
#    paragraph to surround the run in the "结束" paragraph.
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs(16)
$r16 = $d.Range($p16.Range.Start, $p16.Range.End - 1)
$inner16 = '<w:r>' + $rPrHint + '<w:t>访问</w:t></w:r>' +
           '<w:r>' + $rPrHint + '<w:t>https://localhost:8443/</w:t></w:r>'
Insert-Fragment $r16 $inner16

$p17 = $d.Paragraphs(17)
$r17 = $d.Range($p17.Range.Start, $p17.Range.End - 1)
$inner17 = '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
           '<w:r>' + $rPrHint + '<w:t>结束</w:t></w:r>' +
           '<w:bookmarkEnd w:id="0"/>'
Insert-Fragment $r17 $inner17
